# Fruta / hortaliza, semanal
#
# A new weekly observation is inserted as a new row 90 (pushing the
# existing rows 90-128 down to 91-129); the rest of the sheet is
# unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 90 - this shifts rows
# 90:128 down to 91:129 and keeps the surrounding formatting.
$ws.Rows('90:90').Insert()

# Populate the newly inserted row with the new observation.
$ws.Range('A90').Value = 6
$ws.Range('B90').Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range('C90').Value = 'Metropolitana'
$ws.Range('D90').Value = 44455
$ws.Range('E90').Value = 13
$ws.Range('F90').Value = 100112022
$ws.Range('G90').Value = 'Arveja Verde'
$ws.Range('H90').Value = 'Sin especificar'
$ws.Range('I90').Value = 'Primera'
$ws.Range('J90').Value = 400
$ws.Range('K90').Value = 25000
$ws.Range('L90').Value = 27000
$ws.Range('M90').Value = 25850
$ws.Range('N90').Value = '$/malla 25 kilos'
$ws.Range('O90').Value = 'Provincia de Huasco'
$ws.Range('P90').Value = 1034
$ws.Range('Q90').Value = 25
$ws.Range('R90').Value = 'Hortaliza'
